# Updates NATMI LR-pair output (Col5a3-Sdc3) with new TPM-based values:
# a new "Resolving-Mac" sending/target cluster is introduced alongside the
# existing ECs/FAPs/MuSCs clusters, which changes every specificity/weight
# statistic in rows 2-13 (now computed against 4 clusters instead of 3) and
# adds 4 new rows (14-17) for Resolving-Mac as a sending cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Col5a3"
$ws.Cells.Item(2, 3).Value = "Sdc3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.8801540000000001
$ws.Cells.Item(2, 8).Value = 2.640462
$ws.Cells.Item(2, 9).Value = 0.007623830312051026
$ws.Cells.Item(2, 10).Value = 0.007623830312051028
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 29.47402433333333
$ws.Cells.Item(2, 14).Value = 88.422073
$ws.Cells.Item(2, 15).Value = 0.295877356230023
$ws.Cells.Item(2, 16).Value = 0.295877356230023
$ws.Cells.Item(2, 17).Value = 25.94168041308067
$ws.Cells.Item(2, 18).Value = 233.475123717726
$ws.Cells.Item(2, 19).Value = 0.002255718757075969
$ws.Cells.Item(2, 20).Value = 0.002255718757075969

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Col5a3"
$ws.Cells.Item(3, 3).Value = "Sdc3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.8801540000000001
$ws.Cells.Item(3, 8).Value = 2.640462
$ws.Cells.Item(3, 9).Value = 0.007623830312051026
$ws.Cells.Item(3, 10).Value = 0.007623830312051028
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 18.11074233333333
$ws.Cells.Item(3, 14).Value = 54.332227
$ws.Cells.Item(3, 15).Value = 0.1818061388681701
$ws.Cells.Item(3, 16).Value = 0.1818061388681701
$ws.Cells.Item(3, 17).Value = 15.94024230765267
$ws.Cells.Item(3, 18).Value = 143.462180768874
$ws.Cells.Item(3, 19).Value = 0.001386059152420114
$ws.Cells.Item(3, 20).Value = 0.001386059152420114

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Col5a3"
$ws.Cells.Item(4, 3).Value = "Sdc3"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.8801540000000001
$ws.Cells.Item(4, 8).Value = 2.640462
$ws.Cells.Item(4, 9).Value = 0.007623830312051026
$ws.Cells.Item(4, 10).Value = 0.007623830312051028
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 7.238098333333333
$ws.Cells.Item(4, 14).Value = 21.714295
$ws.Cells.Item(4, 15).Value = 0.07266023040422054
$ws.Cells.Item(4, 16).Value = 0.07266023040422054
$ws.Cells.Item(4, 17).Value = 6.370641200476667
$ws.Cells.Item(4, 18).Value = 57.33577080429001
$ws.Cells.Item(4, 19).Value = 0.0005539492670363082
$ws.Cells.Item(4, 20).Value = 0.0005539492670363083

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Col5a3"
$ws.Cells.Item(5, 3).Value = "Sdc3"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.8801540000000001
$ws.Cells.Item(5, 8).Value = 2.640462
$ws.Cells.Item(5, 9).Value = 0.007623830312051026
$ws.Cells.Item(5, 10).Value = 0.007623830312051028
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 44.79281599999999
$ws.Cells.Item(5, 14).Value = 134.378448
$ws.Cells.Item(5, 15).Value = 0.4496562744975863
$ws.Cells.Item(5, 16).Value = 0.4496562744975864
$ws.Cells.Item(5, 17).Value = 39.424576173664
$ws.Cells.Item(5, 18).Value = 354.821185562976
$ws.Cells.Item(5, 19).Value = 0.003428103135518635
$ws.Cells.Item(5, 20).Value = 0.003428103135518637

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Col5a3"
$ws.Cells.Item(6, 3).Value = "Sdc3"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 92.11319466666667
$ws.Cells.Item(6, 8).Value = 276.339584
$ws.Cells.Item(6, 9).Value = 0.7978778323334215
$ws.Cells.Item(6, 10).Value = 0.7978778323334216
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 29.47402433333333
$ws.Cells.Item(6, 14).Value = 88.422073
$ws.Cells.Item(6, 15).Value = 0.295877356230023
$ws.Cells.Item(6, 16).Value = 0.295877356230023
$ws.Cells.Item(6, 17).Value = 2714.946541026404
$ws.Cells.Item(6, 18).Value = 24434.51886923763
$ws.Cells.Item(6, 19).Value = 0.2360739836253543
$ws.Cells.Item(6, 20).Value = 0.2360739836253543

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Col5a3"
$ws.Cells.Item(7, 3).Value = "Sdc3"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 92.11319466666667
$ws.Cells.Item(7, 8).Value = 276.339584
$ws.Cells.Item(7, 9).Value = 0.7978778323334215
$ws.Cells.Item(7, 10).Value = 0.7978778323334216
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 18.11074233333333
$ws.Cells.Item(7, 14).Value = 54.332227
$ws.Cells.Item(7, 15).Value = 0.1818061388681701
$ws.Cells.Item(7, 16).Value = 0.1818061388681701
$ws.Cells.Item(7, 17).Value = 1668.238334108175
$ws.Cells.Item(7, 18).Value = 15014.14500697357
$ws.Cells.Item(7, 19).Value = 0.1450590879850446
$ws.Cells.Item(7, 20).Value = 0.1450590879850446

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Col5a3"
$ws.Cells.Item(8, 3).Value = "Sdc3"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 92.11319466666667
$ws.Cells.Item(8, 8).Value = 276.339584
$ws.Cells.Item(8, 9).Value = 0.7978778323334215
$ws.Cells.Item(8, 10).Value = 0.7978778323334216
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 7.238098333333333
$ws.Cells.Item(8, 14).Value = 21.714295
$ws.Cells.Item(8, 15).Value = 0.07266023040422054
$ws.Cells.Item(8, 16).Value = 0.07266023040422054
$ws.Cells.Item(8, 17).Value = 666.7243607948089
$ws.Cells.Item(8, 18).Value = 6000.51924715328
$ws.Cells.Item(8, 19).Value = 0.05797398713176645
$ws.Cells.Item(8, 20).Value = 0.05797398713176646

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Col5a3"
$ws.Cells.Item(9, 3).Value = "Sdc3"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 92.11319466666667
$ws.Cells.Item(9, 8).Value = 276.339584
$ws.Cells.Item(9, 9).Value = 0.7978778323334215
$ws.Cells.Item(9, 10).Value = 0.7978778323334216
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 44.79281599999999
$ws.Cells.Item(9, 14).Value = 134.378448
$ws.Cells.Item(9, 15).Value = 0.4496562744975863
$ws.Cells.Item(9, 16).Value = 0.4496562744975864
$ws.Cells.Item(9, 17).Value = 4126.009379876181
$ws.Cells.Item(9, 18).Value = 37134.08441888563
$ws.Cells.Item(9, 19).Value = 0.3587707735912561
$ws.Cells.Item(9, 20).Value = 0.3587707735912562

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Col5a3"
$ws.Cells.Item(10, 3).Value = "Sdc3"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 22.42634833333333
$ws.Cells.Item(10, 8).Value = 67.279045
$ws.Cells.Item(10, 9).Value = 0.1942554078175884
$ws.Cells.Item(10, 10).Value = 0.1942554078175884
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 29.47402433333333
$ws.Cells.Item(10, 14).Value = 88.422073
$ws.Cells.Item(10, 15).Value = 0.295877356230023
$ws.Cells.Item(10, 16).Value = 0.295877356230023
$ws.Cells.Item(10, 17).Value = 660.9947364844761
$ws.Cells.Item(10, 18).Value = 5948.952628360285
$ws.Cells.Item(10, 19).Value = 0.05747577649845299
$ws.Cells.Item(10, 20).Value = 0.05747577649845299

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Col5a3"
$ws.Cells.Item(11, 3).Value = "Sdc3"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 22.42634833333333
$ws.Cells.Item(11, 8).Value = 67.279045
$ws.Cells.Item(11, 9).Value = 0.1942554078175884
$ws.Cells.Item(11, 10).Value = 0.1942554078175884
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 18.11074233333333
$ws.Cells.Item(11, 14).Value = 54.332227
$ws.Cells.Item(11, 15).Value = 0.1818061388681701
$ws.Cells.Item(11, 16).Value = 0.1818061388681701
$ws.Cells.Item(11, 17).Value = 406.1578161425795
$ws.Cells.Item(11, 18).Value = 3655.420345283215
$ws.Cells.Item(11, 19).Value = 0.03531682564957749
$ws.Cells.Item(11, 20).Value = 0.03531682564957749

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Col5a3"
$ws.Cells.Item(12, 3).Value = "Sdc3"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 22.42634833333333
$ws.Cells.Item(12, 8).Value = 67.279045
$ws.Cells.Item(12, 9).Value = 0.1942554078175884
$ws.Cells.Item(12, 10).Value = 0.1942554078175884
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 7.238098333333333
$ws.Cells.Item(12, 14).Value = 21.714295
$ws.Cells.Item(12, 15).Value = 0.07266023040422054
$ws.Cells.Item(12, 16).Value = 0.07266023040422054
$ws.Cells.Item(12, 17).Value = 162.3241144942528
$ws.Cells.Item(12, 18).Value = 1460.917030448275
$ws.Cells.Item(12, 19).Value = 0.0141146426892918
$ws.Cells.Item(12, 20).Value = 0.0141146426892918

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Col5a3"
$ws.Cells.Item(13, 3).Value = "Sdc3"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 22.42634833333333
$ws.Cells.Item(13, 8).Value = 67.279045
$ws.Cells.Item(13, 9).Value = 0.1942554078175884
$ws.Cells.Item(13, 10).Value = 0.1942554078175884
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 44.79281599999999
$ws.Cells.Item(13, 14).Value = 134.378448
$ws.Cells.Item(13, 15).Value = 0.4496562744975863
$ws.Cells.Item(13, 16).Value = 0.4496562744975864
$ws.Cells.Item(13, 17).Value = 1004.539294446907
$ws.Cells.Item(13, 18).Value = 9040.85365002216
$ws.Cells.Item(13, 19).Value = 0.0873481629802661
$ws.Cells.Item(13, 20).Value = 0.08734816298026611

# Row 14
$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Col5a3"
$ws.Cells.Item(14, 3).Value = "Sdc3"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.02804566666666667
$ws.Cells.Item(14, 8).Value = 0.084137
$ws.Cells.Item(14, 9).Value = 0.0002429295369390043
$ws.Cells.Item(14, 10).Value = 0.0002429295369390043
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 29.47402433333333
$ws.Cells.Item(14, 14).Value = 88.422073
$ws.Cells.Item(14, 15).Value = 0.295877356230023
$ws.Cells.Item(14, 16).Value = 0.295877356230023
$ws.Cells.Item(14, 17).Value = 0.8266186617778889
$ws.Cells.Item(14, 18).Value = 7.439567956001
$ws.Cells.Item(14, 19).Value = 0.0000718773491396963
$ws.Cells.Item(14, 20).Value = 0.00007187734913969631

# Row 15
$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Col5a3"
$ws.Cells.Item(15, 3).Value = "Sdc3"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.02804566666666667
$ws.Cells.Item(15, 8).Value = 0.084137
$ws.Cells.Item(15, 9).Value = 0.0002429295369390043
$ws.Cells.Item(15, 10).Value = 0.0002429295369390043
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 18.11074233333333
$ws.Cells.Item(15, 14).Value = 54.332227
$ws.Cells.Item(15, 15).Value = 0.1818061388681701
$ws.Cells.Item(15, 16).Value = 0.1818061388681701
$ws.Cells.Item(15, 17).Value = 0.5079278425665555
$ws.Cells.Item(15, 18).Value = 4.571350583099001
$ws.Cells.Item(15, 19).Value = 0.00004416608112791288
$ws.Cells.Item(15, 20).Value = 0.00004416608112791289

# Row 16
$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Col5a3"
$ws.Cells.Item(16, 3).Value = "Sdc3"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.02804566666666667
$ws.Cells.Item(16, 8).Value = 0.084137
$ws.Cells.Item(16, 9).Value = 0.0002429295369390043
$ws.Cells.Item(16, 10).Value = 0.0002429295369390043
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 7.238098333333333
$ws.Cells.Item(16, 14).Value = 21.714295
$ws.Cells.Item(16, 15).Value = 0.07266023040422054
$ws.Cells.Item(16, 16).Value = 0.07266023040422054
$ws.Cells.Item(16, 17).Value = 0.2029972931572222
$ws.Cells.Item(16, 18).Value = 1.826975638415
$ws.Cells.Item(16, 19).Value = 0.00001765131612597866
$ws.Cells.Item(16, 20).Value = 0.00001765131612597866

# Row 17
$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "Col5a3"
$ws.Cells.Item(17, 3).Value = "Sdc3"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.02804566666666667
$ws.Cells.Item(17, 8).Value = 0.084137
$ws.Cells.Item(17, 9).Value = 0.0002429295369390043
$ws.Cells.Item(17, 10).Value = 0.0002429295369390043
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 44.79281599999999
$ws.Cells.Item(17, 14).Value = 134.378448
$ws.Cells.Item(17, 15).Value = 0.4496562744975863
$ws.Cells.Item(17, 16).Value = 0.4496562744975864
$ws.Cells.Item(17, 17).Value = 1.256244386597333
$ws.Cells.Item(17, 18).Value = 11.306199479376
$ws.Cells.Item(17, 19).Value = 0.0001092347905454164
$ws.Cells.Item(17, 20).Value = 0.0001092347905454165
